# Updated cryptos list on Wed Jul 31 04:54:39 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" values in column D are plain numeric-looking strings (e.g. "582.48").
# Excel's COM Value setter auto-converts such strings to real numbers, which would
# change the cell type / introduce float rounding noise. To keep them as literal
# text (matching the source data, which stores prices as text), we mark those
# specific cells as Text ("@") before writing the value. Values that contain two
# "." separators (e.g. "65.610.34") are never auto-parsed as numbers, so they
# don't need this treatment.
$textPriceCells = @("D5","D6","D7","D11","D15","D16","D18","D19","D21","D22","D23","D24","D26","D27","D29","D31","D32","D34","D35","D37","D40","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "65.610.34"
$ws.Range("E2").Value = "  -1.58%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.270.84"
$ws.Range("E3").Value = "  -1.21%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.04%  "

# Row 5 - BNB
$ws.Range("D5").Value = "582.48"
$ws.Range("E5").Value = "  +1.76%  "

# Row 6 - Solana
$ws.Range("D6").Value = "178.88"
$ws.Range("E6").Value = "  -1.94%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.642"
$ws.Range("E7").Value = "  +7.12%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.01%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -3.95%  "

# Row 10 - Toncoin
$ws.Range("E10").Value = "  +1.39%  "

# Row 11 - Cardano
$ws.Range("D11").Value = "0.401"
$ws.Range("E11").Value = "  -0.66%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "3.839.50"
$ws.Range("E12").Value = "  -1.17%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  -4.52%  "

# Row 14 - WrappedBTC
$ws.Range("D14").Value = "65.748.21"
$ws.Range("E14").Value = "  -1.41%  "

# Row 15 - Avalanche
$ws.Range("D15").Value = "25.88"
$ws.Range("E15").Value = "  -4.67%  "

# Row 16 and 17 swap: WrappedEther <-> ShibaInu
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "0.0000162"
$ws.Range("E16").Value = "  -3.21%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.257.30"
$ws.Range("E17").Value = "  -1.86%  "

# Row 18 - BitcoinCash
$ws.Range("D18").Value = "426.27"
$ws.Range("E18").Value = "  -1.59%  "

# Row 19 - Chainlink
$ws.Range("D19").Value = "13.16"
$ws.Range("E19").Value = "  -4.15%  "

# Row 20 - Polkadot
$ws.Range("E20").Value = "  -3.71%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "7.34"
$ws.Range("E21").Value = "  -3.71%  "

# Row 22 - Litecoin
$ws.Range("D22").Value = "71.84"
$ws.Range("E22").Value = "  -2.59%  "

# Row 23 - Dai
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.07%  "

# Row 24 - LEO
$ws.Range("D24").Value = "5.68"
$ws.Range("E24").Value = "  +0.20%  "

# Row 25 - WrappedeETH
$ws.Range("D25").Value = "3.428.72"
$ws.Range("E25").Value = "  -0.86%  "

# Row 26 - Polygon
$ws.Range("D26").Value = "0.507"
$ws.Range("E26").Value = "  -1.64%  "

# Row 27 - Kaspa
$ws.Range("D27").Value = "0.196"
$ws.Range("E27").Value = "  +0.82%  "

# Row 28 - PEPE
$ws.Range("E28").Value = "  -5.34%  "

# Row 29 - InternetComputer(DFINITY)
$ws.Range("D29").Value = "8.82"
$ws.Range("E29").Value = "  -2.61%  "

# Row 30 - Binance-PegBSC-USD
$ws.Range("E30").Value = "  +0.03%  "

# Row 31 - PancakeSwap
$ws.Range("D31").Value = "1.95"
$ws.Range("E31").Value = "  -0.47%  "

# Row 32 - EthereumClassic
$ws.Range("D32").Value = "22.14"
$ws.Range("E32").Value = "  -2.80%  "

# Row 33 - USDe
$ws.Range("E33").Value = "  -0.01%  "

# Row 34 - NEARProtocol
$ws.Range("D34").Value = "5.12"
$ws.Range("E34").Value = "  -4.01%  "

# Row 35 - Aptos
$ws.Range("D35").Value = "6.56"
$ws.Range("E35").Value = "  -3.13%  "

# Row 36 - Fetch.AI
$ws.Range("E36").Value = "  -3.74%  "

# Row 37 - Monero
$ws.Range("D37").Value = "159.36"
$ws.Range("E37").Value = "  -0.59%  "

# Row 38 - ImmutableX (no changes)

# Row 39 - Stacks
$ws.Range("E39").Value = "  -3.29%  "

# Row 40 - EnergySwap
$ws.Range("D40").Value = "26.23"
$ws.Range("E40").Value = "  -3.84%  "

# Row 41 - Maker
$ws.Range("D41").Value = "2.781.10"
$ws.Range("E41").Value = "  -1.08%  "

# Row 42 - Mantle
$ws.Range("E42").Value = "  -3.22%  "

# Row 43 - Filecoin
$ws.Range("D43").Value = "4.30"
$ws.Range("E43").Value = "  -3.27%  "

# Row 44 - OKB
$ws.Range("D44").Value = "40.02"
$ws.Range("E44").Value = "  -0.28%  "

# Row 45 - Hedera
$ws.Range("D45").Value = "0.0655"
$ws.Range("E45").Value = "  -2.94%  "

# Row 46 - RenderToken
$ws.Range("D46").Value = "5.90"
$ws.Range("E46").Value = "  -5.51%  "

# Row 47 - dogwifhat
$ws.Range("D47").Value = "2.27"
$ws.Range("E47").Value = "  -2.99%  "

# Row 48 - Bittensor
$ws.Range("D48").Value = "314.15"
$ws.Range("E48").Value = "  -1.79%  "

# Row 49 - InjectiveProtocol
$ws.Range("D49").Value = "22.96"
$ws.Range("E49").Value = "  -5.73%  "

# Row 50 - VeChain
$ws.Range("D50").Value = "0.0266"
$ws.Range("E50").Value = "  -2.59%  "

# Row 51 - Stellar
$ws.Range("D51").Value = "0.104"
$ws.Range("E51").Value = "  +4.76%  "
